$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New helper rows 38-45 on MaxDed1 (sheet2) ---
# Labels (column B) are entered in this specific order so that the shared-string
# table is populated in the same first-use order as the authored workbook:
#   33 x.loss                  -> B38
#   34 deductible               -> B41
#   35 loss                     -> B42
#   36 accumulated_limit        -> B43
#   37 x.limit_surplus          -> B39
#   38 x.effective_deductible   -> B40
#   39 adjusted limit_surplus   -> B45
#   40 loss_adj                 -> B44
$ws.Range("B38").Value = "x.loss"
$ws.Range("B41").Value = "deductible"
$ws.Range("B42").Value = "loss"
$ws.Range("B43").Value = "accumulated_limit"
$ws.Range("B39").Value = "x.limit_surplus"
$ws.Range("B40").Value = "x.effective_deductible"
$ws.Range("B45").Value = "adjusted limit_surplus"
$ws.Range("B44").Value = "loss_adj"

# Formulas (column C), entered in row order. Each gets the thousands-separator
# number format (matches the new cellXfs entry, numFmtId 3, used throughout
# this helper block).
$ws.Range("C38").Formula = "=F28"
$ws.Range("C38").NumberFormat = "#,##0"
$ws.Range("C39").Formula = "=F23-F26"
$ws.Range("C39").NumberFormat = "#,##0"
$ws.Range("C40").Formula = "=F24"
$ws.Range("C40").NumberFormat = "#,##0"
$ws.Range("C41").Formula = "=C16"
$ws.Range("C41").NumberFormat = "#,##0"
$ws.Range("C42").Formula = "=C38+C40-C41"
$ws.Range("C42").NumberFormat = "#,##0"
$ws.Range("C43").Formula = "=IF(C39>0,C38)"
$ws.Range("C43").NumberFormat = "#,##0"
$ws.Range("C44").Formula = "=C42-C38"
$ws.Range("C44").NumberFormat = "#,##0"
$ws.Range("C45").Formula = "=C39+C44"
$ws.Range("C45").NumberFormat = "#,##0"

# Stray formatted-but-empty cells that accompany the new block.
$ws.Range("G28").NumberFormat = "#,##0"
$ws.Range("G29").NumberFormat = "#,##0"
$ws.Range("F44").NumberFormat = "#,##0"
$ws.Range("F45").NumberFormat = "#,##0"

# Tidy up column widths for the newly-used columns C and G.
$ws.Columns.Item(3).AutoFit() | Out-Null
$ws.Columns.Item(7).AutoFit() | Out-Null

# Restore the view: selection on C45, scrolled so row 19 is at the top.
$ws.Range("C45").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
